# GL07 Postback.xlsx - "Add files via upload" edit
#
# The author cleared out the example/demo values that had been entered into
# the Controls sheet's yellow input cells (batch id, client, period, the
# TODAY() voucher/transaction date, voucher type, plus the two spare rows),
# leaving the sheet blank/ready for the next user. Everything else in the
# workbook (the GL07 / GL07 (2) report sheets, shared strings, etc.) just
# follows on from that - their cached formula values and shared-string
# indices shift automatically once the source cells are cleared.
#
# The commit also carries forward some Excel UI/view bookkeeping: the last
# selected cell on each sheet, and which sheet/tab was active when it was
# saved (GL07 (2) instead of GL07).

$wb = $excel.ActiveWorkbook

# --- Controls sheet: clear the example values entered in the input cells ---
$controls = $wb.Worksheets.Item("Controls")
$controls.Activate()

$controls.Range("C5").ClearContents()   # Client            (was "FE")
$controls.Range("C6").ClearContents()   # Batch ID          (was "DE2")
$controls.Range("C7").ClearContents()   # Period            (was 201401)
$controls.Range("C8").ClearContents()   # Trans/Value Date  (was =TODAY())
$controls.Range("C9").ClearContents()   # Voucher Type      (was "GA")
$controls.Range("C10").ClearContents()  # Currency          (was "GBP")
$controls.Range("C11").ClearContents()  # (spare row)       (was "BA")
$controls.Range("C12").ClearContents()  # Voucher Date      (was =C8)

$controls.Range("X14").Select()

# --- GL07 sheet: just update the remembered selection / scroll position ---
$gl07 = $wb.Worksheets.Item("GL07")
$gl07.Activate()
$gl07.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

# --- GL07 (2) sheet: becomes the active/selected tab on save ---
$gl072 = $wb.Worksheets.Item("GL07 (2)")
$gl072.Activate()
$gl072.Range("F29").Select()
